$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "F-SW-SD-10"
$oldName = "Sheet1"
$newName = "F-SW-SD-10"
$ws = $wb.Worksheets.Item($oldName)
$ws.Name = $newName

# The workbook-level Print_Area defined name still points at the old sheet
# name after a plain rename, so repoint it explicitly at the renamed sheet.
$printArea = $wb.Names.Item("$newName!Print_Area")
$printArea.RefersTo = "='$newName'!`$A`$1:`$G`$22"
